# Generate Report for Handoff
# Updates the Priority and Latest Handoff Datetime for the four files that were
# picked up for handoff (5e9efd1d..., 62f5b9fe..., 693c19e0..., 8472a05a...)
# on both the zh-cn and de-de localization status sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# zh-cn sheet: rows 4-7 -> Priority (E) moves from "low" to "ht",
# and Latest Handoff Datetime (H) is refreshed to the new handoff time.
foreach ($row in 4..7) {
    $wsZhCn.Range("E$row").Value = "ht"
    $wsZhCn.Range("H$row").Value = "2016-09-02 08:35:55"
}

# de-de sheet: rows 4-7 -> Priority (E) moves from "low" to "ht".
# Latest Handoff Datetime (H) here shares the other handoff timestamp
# string, refreshed below for the whole workbook.
foreach ($row in 4..7) {
    $wsDeDe.Range("E$row").Value = "ht"
}

# The "Latest HO Xliff Generate Date" / de-de handoff timestamp text itself is
# refreshed from 2016-09-02 08:35:44 to 2016-09-02 08:35:59.
foreach ($row in 4..7) {
    $wsDeDe.Range("H$row").Value = "2016-09-02 08:35:59"
}

$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($row in 4..7) {
    $wsOverview.Range("G$row").Value = "2016-09-02 08:35:59"
}
